{"js": "// Citation-check update: replace the placeholder \"Ref-XXXXXX\" style in-text\n// citation markers with their resolved bibliographic citations (or, for two\n// references that stay as generated IDs, with their updated ID values).\n//\n// Because the same placeholder text \"Ref-J7X2BZ\" occurs twice in the\n// document (paragraph 1 and paragraph 2) but resolves to two different\n// replacement strings, a single document-wide search/replace would be\n// wrong. Instead we scope each search to the specific paragraph that needs\n// to change, then replace just that match.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst replacements = [\n  { paraIndex: 0, find: \"Ref-J7X2BZ\", replace: \"Brown & Garcia, 2018\" },\n  { paraIndex: 1, find: \"Ref-J7X2BZ\", replace: \"Ref-u170605\" },\n  { paraIndex: 2, find: \"Ref-DJ74KL\", replace: \"Ref-u746170\" },\n  { paraIndex: 3, find: \"Ref-DJ49F2\", replace: \"Ref-u805438\" },\n  { paraIndex: 4, find: \"Ref-J7X8A2\", replace: \"Pearse et al., 2001\" },\n];\n\nfor (const r of replacements) {\n  const para = paragraphs.items[r.paraIndex];\n  const results = para.search(r.find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  results.items[0].insertText(r.replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Citation-check update: replace the placeholder \"Ref-XXXXXX\" style in-text\n# citation markers with their resolved bibliographic citations (or, for two\n# references that stay as generated IDs, with their updated ID values).\n#\n# Because the same placeholder text \"Ref-J7X2BZ\" occurs twice in the\n# document (paragraph 1 and paragraph 2) but resolves to two different\n# replacement strings, a single document-wide Find/Replace would be wrong.\n# Instead we scope each Find/Replace to the specific paragraph that needs\n# to change.\n\n$d = $word.ActiveDocument\n\nfunction Replace-InParagraph($paraIndex, $findText, $replaceText) {\n    $rng = $d.Paragraphs($paraIndex).Range\n    $find = $rng.Find\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute(\n        [ref]$find.Text,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        [ref]$find.Replacement.Text,\n        2\n    ) | Out-Null\n}\n\nReplace-InParagraph 1 \"Ref-J7X2BZ\" \"Brown & Garcia, 2018\"\nReplace-InParagraph 2 \"Ref-J7X2BZ\" \"Ref-u170605\"\nReplace-InParagraph 3 \"Ref-DJ74KL\" \"Ref-u746170\"\nReplace-InParagraph 4 \"Ref-DJ49F2\" \"Ref-u805438\"\nReplace-InParagraph 5 \"Ref-J7X8A2\" \"Pearse et al., 2001\"\n"}
